{"js": "// Day 2 writing automatic tests using Behave\n// Highlight the \"Welcome screen...\" test-scenario bullet in green.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) =>\n    p.text.indexOf(\n      \"Welcome screen should be clear, readable, visible on the page when finishes loading.\"\n    ) !== -1\n);\n\nif (!target) {\n  throw new Error(\n    \"Could not find the 'Welcome screen...' test-scenario paragraph.\"\n  );\n}\n\n// Set the highlight through the paragraph's Font so that both the run and\n// the paragraph mark's run-properties (pPr/rPr) pick up the highlight \u2014\n// matching what Word does when you select the whole paragraph (pilcrow\n// included) and apply a highlighter color.\ntarget.font.highlightColor = \"BrightGreen\";\n\nawait context.sync();\n", "ps1": "# Day 2 writing automatic tests using Behave\n# Highlight the \"Welcome screen...\" test-scenario bullet in green, the\n# same way Word does when you select the paragraph (including its\n# paragraph mark) and apply a highlight color from the ribbon.\n\n$d = $word.ActiveDocument\n\n$target = $d.Paragraphs | Where-Object {\n    $_.Range.Text -like \"*Welcome screen should be clear, readable, visible on the page when finishes loading.*\"\n}\n\nif (-not $target) {\n    throw \"Could not find the 'Welcome screen...' test-scenario paragraph.\"\n}\n\n# Apply the highlight to the whole paragraph range (Paragraph.Range already\n# spans the text plus its trailing paragraph mark), via Font so that both\n# the run and the paragraph mark's run-properties pick up the highlight.\n$target.Range.Font.HighlightColorIndex = \"wdBrightGreen\"\n"}
